$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.036978960037231
$ws.Range("B1").Value = 6.341531276702881
$ws.Range("C1").Value = 6.914728164672852
$ws.Range("D1").Value = 7.404232025146484
$ws.Range("E1").Value = 4.739113807678223
